# Fruta / hortaliza, semanal
#
# The weekly refresh re-shuffles the per-record fields (Fecha, Volumen,
# Precio minimo/maximo/promedio ponderado, Origen, Precio $/Kg) across the
# existing data rows (2-15) of the active sheet. Row 4 is left untouched.
# Apply the new values exactly as they land after the shuffle.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @{
    2  = @{ D = 44617; M = 90;  N = 6500; O = 6500; P = 6500; R = "Provincia de Curicó";  S = 3250 }
    3  = @{ D = 44231; M = 150; N = 3400; O = 3400; P = 3400; R = "Provincia de Curicó";  S = 1700 }
    5  = @{ D = 44188; M = 150; N = 3000; O = 3400; P = 3240; R = "Provincia de Linares"; S = 1620 }
    6  = @{ D = 44238; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó";  S = 1900 }
    7  = @{ D = 44586; M = 250; N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó";  S = 2500 }
    8  = @{ D = 44533; M = 150; N = 4000; O = 4000; P = 4000; R = "Provincia de Curicó";  S = 2000 }
    9  = @{ D = 44582; M = 380; N = 5000; O = 5000; P = 5000; R = "Provincia de Curicó";  S = 2500 }
    10 = @{ D = 44237; M = 100; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó";  S = 1900 }
    11 = @{ D = 44168; M = 170; N = 8000; O = 8000; P = 8000; R = "Provincia de Linares"; S = 4000 }
    12 = @{ D = 44174; M = 200; N = 3200; O = 3200; P = 3200; R = "Provincia de Curicó";  S = 1600 }
    13 = @{ D = 44208; M = 85;  N = 3000; O = 3000; P = 3000; R = "Provincia de Linares"; S = 1500 }
    14 = @{ D = 44236; M = 300; N = 3600; O = 4000; P = 3800; R = "Provincia de Curicó";  S = 1900 }
    15 = @{ D = 44232; M = 200; N = 3000; O = 3000; P = 3000; R = "Provincia de Curicó";  S = 1500 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("N$r").Value = $vals.N
    $ws.Range("O$r").Value = $vals.O
    $ws.Range("P$r").Value = $vals.P
    $ws.Range("R$r").Value = $vals.R
    $ws.Range("S$r").Value = $vals.S
}
